$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Prix TSP" column to the existing table (Tableau1) ---
$lo = $ws.ListObjects.Item(1)
$newCol = $lo.ListColumns.Add()

# Give the new header cell the same formatting as the other header cells (C1/D1/E1)
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Set header text (also renames the table column automatically)
$ws.Range("F1").Value2 = "Prix TSP"

# Column E gets a custom width in the new file
$ws.Columns.Item(5).ColumnWidth = 15.333333333333334

# --- Fill in the "Prix TSP" values for every data row ---
$tsp = @{
  2  = 267.74737611748355
  3  = 320.50615838606103
  4  = 335.36322723603217
  5  = 411.58927879408168
  6  = 481.93721062003419
  7  = 493.58760761581902
  8  = 538.67783908274714
  9  = 619.81446658599202
  10 = 583.17754284926309
  11 = 547.08557077308285
  12 = 606.38191056603966
  13 = 644.12187851407236
  14 = 701.13260620976212
  15 = 717
  16 = 700.64750792612949
  17 = 903.49661935509459
  18 = 970.96058802397135
  19 = 717.76783814089595
  20 = 721.78502775499192
  21 = 785.74556336099408
  22 = 693.43033122426732
  23 = 661.2415008512379
  24 = 387.16212765957448
  25 = 439.59987941469245
  26 = 473.7031495401406
  27 = 443.04752604764099
  28 = 383.01967167590209
  29 = 384.37074369159416
  30 = 380.52525067592018
  31 = 331.99481299422956
  32 = 342.00988593704562
  33 = 339.80463588897572
  34 = 376.09120144127439
  35 = 396.24926974395703
  36 = 401.7686977749708
  37 = 398.20689444428109
  38 = 406.14347451452056
  39 = 397.24050742441119
  40 = 395.91180800645901
  41 = 358.73110949800457
  42 = 381.68920664576007
  43 = 384.25059914910969
  44 = 380.15482969888251
  45 = 412.59185125703368
  46 = 428.32217072847618
  47 = 431.42029846944047
  48 = 422.64952432346183
  49 = 436.99855262881238
}

foreach ($row in $tsp.Keys) {
  $ws.Cells.Item($row, 6).Value2 = $tsp[$row]
}

# --- Match the final selection recorded in the saved workbook ---
$ws.Range("I45").Select() | Out-Null
